# Indicadores Plano de Trabalho Junho 23
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column J header + formatting -----------------------------------
# Copy the existing bordered style (s=1, used by B4 etc.) onto the header
# row (A3:I3, which previously had no style) and onto column J (rows 3-13)
# so the new cells reuse the same cellXf rather than creating duplicates.
$ws.Range("B4").Copy()
$ws.Range("A3:I3").PasteSpecial(-4122)

$ws.Range("B4").Copy()
$ws.Range("J3:J14").PasteSpecial(-4122)
$ws.Range("J3").WrapText = $true

$excel.CutCopyMode = 0

# Row 3 is now taller to accommodate the wrapped header text.
$ws.Rows("3").RowHeight = 51

# New header text for column J.
$ws.Range("J3").Value = "VMP Incluídos Portal"

# --- Data edits -----------------------------------------------------------
$ws.Range("C7").Value = 2007
$ws.Range("G7").Value = 253
$ws.Range("H7").Value = 612
$ws.Range("I7").Value = 1729
$ws.Range("J7").Value = 136

$ws.Range("J13").Value = " "

# --- Row 14 totals ----------------------------------------------------------
$ws.Range("D14").Formula = "=C14/B14"
$ws.Range("F14").Formula = "=SUM(F4:F13)"
$ws.Range("G14").Formula = "=SUM(G4:G13)"
$ws.Range("H14").Formula = "=SUM(H4:H13)"
$ws.Range("I14").Formula = "=SUM(I4:I13)"
$ws.Range("J14").Formula = "=SUM(J5:J13)"

# --- Selection update -------------------------------------------------------
$ws.Range("D14").Select()

$wb.Save()
